$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the formatting of the
# existing header cells (bold font, thin border, centered/top alignment).
$ws.Range("I1").Value = "I0"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

$ws.Range("J1").Value = "IF"

# Fill in data rows 2..27: column I (I0) and column J (IF)
for ($r = 2; $r -le 27; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    if ($r -eq 26) {
        $ws.Cells.Item($r, 9).Value = 4
        $ws.Cells.Item($r, 10).Value = 9
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $ipValue
    }
}
